$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.319712
$ws.Range("H2").Value = 3.959136
$ws.Range("M2").Value = 1.953983333333333
$ws.Range("N2").Value = 5.86195
$ws.Range("O2").Value = 0.20183677855562
$ws.Range("P2").Value = 0.20183677855562
$ws.Range("Q2").Value = 2.5786952528
$ws.Range("R2").Value = 23.2082572752
$ws.Range("S2").Value = 0.20183677855562
$ws.Range("T2").Value = 0.20183677855562

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.319712
$ws.Range("H3").Value = 3.959136
$ws.Range("O3").Value = 0.02540617157539598
$ws.Range("P3").Value = 0.02540617157539598
$ws.Range("Q3").Value = 0.324592844288
$ws.Range("R3").Value = 2.921335598592
$ws.Range("S3").Value = 0.02540617157539598
$ws.Range("T3").Value = 0.02540617157539598

# Row 4 (Target cluster: MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.319712
$ws.Range("H4").Value = 3.959136
$ws.Range("M4").Value = 7.481066666666666
$ws.Range("N4").Value = 22.4432
$ws.Range("O4").Value = 0.772757049868984
$ws.Range("P4").Value = 0.772757049868984
$ws.Range("Q4").Value = 9.872853452799999
$ws.Range("R4").Value = 88.8556810752
$ws.Range("S4").Value = 0.772757049868984
$ws.Range("T4").Value = 0.772757049868984
